$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the "Förändrad" date (column C) by one day for all data rows (2-36)
for ($r = 2; $r -le 36; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}

# Swap the designation (A) and area (G) values between rows 35 and 36
$a35 = $ws.Cells.Item(35, 1).Value2
$a36 = $ws.Cells.Item(36, 1).Value2
$ws.Cells.Item(35, 1).Value2 = $a36
$ws.Cells.Item(36, 1).Value2 = $a35

$g35 = $ws.Cells.Item(35, 7).Value2
$g36 = $ws.Cells.Item(36, 7).Value2
$ws.Cells.Item(35, 7).Value2 = $g36
$ws.Cells.Item(36, 7).Value2 = $g35
